$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F16").Value = "placeholder"

$ws.Range("F13:F16").Select()
